# This script applies the diff: it turns the empty, numbered (numId=10)
# "ListParagraph" placeholder (the last bullet item, currently with no text)
# into the bold+italic heading "Prime Number of Set Bits in Binary Representation:"
# and appends a C++ code listing below it (as plain, unnumbered paragraphs),
# directly before the two trailing blank paragraphs at the end of the document.

$d = $word.ActiveDocument

# Plain code lines belonging to the new "Prime Number of Set Bits" listing.
$codeLines = @(
    'class Solution {',
    'public:',
    '    bool isprime(int n)',
    '    {',
    '        if(n<2) return false;',
    '        int i;',
    '        for(i=2; i<=n/2; i++)',
    '            if(n%i == 0) return false;',
    '        return true;',
    '    }',
    '    int countPrimeSetBits(int L, int R) {',
    '        int i, bit, count = 0, n;',
    '        for(i=L; i<=R; i++)',
    '        {',
    '            bit = 0;',
    '            n = i;',
    '            while(n)',
    '            {',
    '                if(n & 1 == 1) bit++;',
    '                n = n>>1;',
    '            }',
    '            if(isprime(bit)) count++;',
    '        }',
    '        return count;',
    '    }',
    '};'
)

# Locate the empty numbered placeholder paragraph (style "List Paragraph",
# numId 10) that currently has no text -- it is the last paragraph in the
# document that still carries that numbering.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.ListFormat.ListType -ne 0 -and $para.Range.Text.Trim().Length -eq 0) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the empty numbered placeholder paragraph"
}

# Build the code-listing paragraphs directly after the placeholder. Each new
# paragraph is immediately stripped of the numbering/List Paragraph style it
# inherits from the split, so it ends up as a plain paragraph (matching the
# rest of the code blocks already present in the document).
$idx = $targetIndex
foreach ($line in $codeLines) {
    $cur = $d.Paragraphs.Item($idx)
    $cur.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Range.ListFormat.RemoveNumbers()
    $newPara.Style = "Normal"
    $newPara.Range.Text = $line
}

# Now turn the placeholder paragraph itself into the bold+italic heading.
$headingPara = $d.Paragraphs.Item($targetIndex)
$headingRange = $headingPara.Range
$headingRange.Text = 'Prime Number of Set Bits in Binary Representation' + ":"
$headingRange.Font.Bold = 1
$headingRange.Font.Italic = 1
$headingRange.Font.BoldBi = 1
$headingRange.Font.ItalicBi = 1

# The heading text is authored as two runs (the title, then a trailing
# colon). Re-toggle the formatting on just the final ":" character so it is
# emitted as its own run, matching the original two-run layout.
$headingPara2 = $d.Paragraphs.Item($targetIndex)
$pStart = $headingPara2.Range.Start
$pEnd = $headingPara2.Range.End
$colonStart = $pEnd - 2
$colonRange = $d.Range($colonStart, $colonStart + 1)
$colonRange.Font.Bold = 0
$colonRange.Font.Bold = 1
$colonRange.Font.Italic = 0
$colonRange.Font.Italic = 1
$colonRange.Font.BoldBi = 0
$colonRange.Font.BoldBi = 1
$colonRange.Font.ItalicBi = 0
$colonRange.Font.ItalicBi = 1
